$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.016490782005725
$ws.Range("D2").Value = 1.021944675990463
$ws.Range("E2").Value = 1.044482531898008
$ws.Range("F2").Value = 1.047566473697818
$ws.Range("I2").Value = 1.026394172581937
$ws.Range("J2").Value = 1.021710291576711
$ws.Range("K2").Value = 1.024780511540088
$ws.Range("L2").Value = 1.047253457019663
$ws.Range("M2").Value = 1.050328752574698
$ws.Range("N2").Value = 1.011365970100988
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.017325744020829
$ws.Range("D3").Value = 1.022532671439952
$ws.Range("E3").Value = 1.045707572248977
$ws.Range("F3").Value = 1.048853260925318
$ws.Range("I3").Value = 1.026495617710091
$ws.Range("J3").Value = 1.022181330669626
$ws.Range("K3").Value = 1.025175665944069
$ws.Range("L3").Value = 1.048288665045977
$ws.Range("M3").Value = 1.051426179047708
$ws.Range("N3").Value = 1.011520875876483
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.017865832929648
$ws.Range("D4").Value = 1.02291245255349
$ws.Range("E4").Value = 1.046501006417348
$ws.Range("F4").Value = 1.049686584711606
$ws.Range("I4").Value = 1.026559124306531
$ws.Range("J4").Value = 1.022485311882561
$ws.Range("K4").Value = 1.025429971193584
$ws.Range("L4").Value = 1.048958702022266
$ws.Range("M4").Value = 1.052136423150934
$ws.Range("N4").Value = 1.011620837376528
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.018092840402293
$ws.Range("D5").Value = 1.023071945458091
$ws.Range("E5").Value = 1.046834747017682
$ws.Range("F5").Value = 1.050037079389938
$ws.Range("I5").Value = 1.026585310463632
$ws.Range("J5").Value = 1.022612910171557
$ws.Range("K5").Value = 1.025536548127677
$ws.Range("L5").Value = 1.049240431483279
$ws.Range("M5").Value = 1.052435043057916
$ws.Range("N5").Value = 1.011662795554914
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.018130953233208
$ws.Range("D6").Value = 1.023098715174785
$ws.Range("E6").Value = 1.046890794170288
$ws.Range("F6").Value = 1.050095938715464
$ws.Range("I6").Value = 1.026589677184912
$ws.Range("J6").Value = 1.022634323000885
$ws.Range("K6").Value = 1.025554423302937
$ws.Range("L6").Value = 1.049287737817372
$ws.Range("M6").Value = 1.052485184628842
$ws.Range("N6").Value = 1.011669836659839
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.017868866396319
$ws.Range("D7").Value = 1.022914584364006
$ws.Range("E7").Value = 1.046505465163545
$ws.Range("F7").Value = 1.049691267387613
$ws.Range("I7").Value = 1.026559476220203
$ws.Range("J7").Value = 1.022487017626107
$ws.Range("K7").Value = 1.025431396589718
$ws.Range("L7").Value = 1.048962466323602
$ws.Range("M7").Value = 1.052140413192985
$ws.Range("N7").Value = 1.011621398282055
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.016772999488573
$ws.Range("D8").Value = 1.022143533986159
$ws.Range("E8").Value = 1.044896385121828
$ws.Range("F8").Value = 1.048001207961682
$ws.Range("I8").Value = 1.026428898034902
$ws.Range("J8").Value = 1.021869649166498
$ws.Range("K8").Value = 1.024914342120365
$ws.Range("L8").Value = 1.047603271981963
$ws.Range("M8").Value = 1.050699605339965
$ws.Range("N8").Value = 1.011418377555503
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.014840553048117
$ws.Range("D9").Value = 1.020779619826178
$ws.Range("E9").Value = 1.042066686185896
$ws.Range("F9").Value = 1.045028317909731
$ws.Range("I9").Value = 1.026182486586354
$ws.Range("J9").Value = 1.020775582027314
$ws.Range("K9").Value = 1.023992661175833
$ws.Range("L9").Value = 1.045209599833458
$ws.Range("M9").Value = 1.048161722719234
$ws.Range("N9").Value = 1.011058552323739
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.013551381450126
$ws.Range("D10").Value = 1.019866922937086
$ws.Range("E10").Value = 1.040183988080828
$ws.Range("F10").Value = 1.043049826915368
$ws.Range("I10").Value = 1.026007291232504
$ws.Range("J10").Value = 1.020042095733855
$ws.Range("K10").Value = 1.023371171030315
$ws.Range("L10").Value = 1.043614704063343
$ws.Range("M10").Value = 1.046470424950998
$ws.Range("N10").Value = 1.010817289971633
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.012992960600715
$ws.Range("D11").Value = 1.019470922759566
$ws.Range("E11").Value = 1.039369636755735
$ws.Range("F11").Value = 1.042193917499635
$ws.Range("I11").Value = 1.02592884939453
$ws.Range("J11").Value = 1.019723523811519
$ws.Range("K11").Value = 1.023100403763114
$ws.Range("L11").Value = 1.042924293292856
$ws.Range("M11").Value = 1.045738210033947
$ws.Range("N11").Value = 1.010712497194097
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.012785508920912
$ws.Range("D12").Value = 1.019323712499149
$ws.Range("E12").Value = 1.039067279566891
$ws.Range("F12").Value = 1.041876111891497
$ws.Range("I12").Value = 1.025899325580152
$ws.Range("J12").Value = 1.019605047419215
$ws.Range("K12").Value = 1.022999580708523
$ws.Range("L12").Value = 1.042667871242613
$ws.Range("M12").Value = 1.045466250908829
$ws.Range("N12").Value = 1.010673523991887
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.012830009322662
$ws.Range("D13").Value = 1.019355294923092
$ws.Range("E13").Value = 1.039132130389322
$ws.Range("F13").Value = 1.04194427699062
$ws.Range("I13").Value = 1.025905676040108
$ws.Range("J13").Value = 1.019630467538652
$ws.Range("K13").Value = 1.023021218803382
$ws.Range("L13").Value = 1.042722873429563
$ws.Range("M13").Value = 1.045524586264056
$ws.Range("N13").Value = 1.010681886066473
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.012975813172986
$ws.Range("D14").Value = 1.019458756715446
$ws.Range("E14").Value = 1.039344641189966
$ws.Range("F14").Value = 1.042167645205373
$ws.Range("I14").Value = 1.025926416837633
$ws.Range("J14").Value = 1.019713733461305
$ws.Range("K14").Value = 1.023092074751791
$ws.Range("L14").Value = 1.042903096817447
$ws.Range("M14").Value = 1.045715729445363
$ws.Range("N14").Value = 1.010709276644871
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.01306564387964
$ws.Range("D15").Value = 1.019522487326953
$ws.Range("E15").Value = 1.039475593178215
$ws.Range("F15").Value = 1.042305285228566
$ws.Range("I15").Value = 1.025939144665688
$ws.Range("J15").Value = 1.01976501720668
$ws.Range("K15").Value = 1.023135698604227
$ws.Range("L15").Value = 1.043014141968784
$ws.Range("M15").Value = 1.045833501431213
$ws.Range("N15").Value = 1.010726146465604
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.013588437843347
$ws.Range("D16").Value = 1.019893187512726
$ws.Range("E16").Value = 1.040238052088251
$ws.Range("F16").Value = 1.043106647328695
$ws.Range("I16").Value = 1.026012442843095
$ws.Range("J16").Value = 1.020063217983952
$ws.Range("K16").Value = 1.023389106121243
$ws.Range("L16").Value = 1.043660528241778
$ws.Range("M16").Value = 1.046519022245058
$ws.Range("N16").Value = 1.010824237907312
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.013916319492195
$ws.Range("D17").Value = 1.020125505781024
$ws.Range("E17").Value = 1.040716554143817
$ws.Range("F17").Value = 1.043609531137336
$ws.Range("I17").Value = 1.026057730358526
$ws.Range("J17").Value = 1.020250012962366
$ws.Range("K17").Value = 1.023547618831042
$ws.Range("L17").Value = 1.044066039166782
$ws.Range("M17").Value = 1.04694906483835
$ws.Range("N17").Value = 1.010885681365246
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.014107547853818
$ws.Range("D18").Value = 1.020260936156222
$ws.Range("E18").Value = 1.040995740277609
$ws.Range("F18").Value = 1.043902931099921
$ws.Range("I18").Value = 1.026083896711704
$ws.Range("J18").Value = 1.020358873910911
$ws.Range("K18").Value = 1.023639916576588
$ws.Range("L18").Value = 1.044302585207453
$ws.Range("M18").Value = 1.047199913852345
$ws.Range("N18").Value = 1.010921488959045
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.014172748455903
$ws.Range("D19").Value = 1.020307101318464
$ws.Range("E19").Value = 1.041090949823438
$ws.Range("F19").Value = 1.044002985899102
$ws.Range("I19").Value = 1.026092776489063
$ws.Range("J19").Value = 1.020395976815218
$ws.Range("K19").Value = 1.023671360545586
$ws.Range("L19").Value = 1.044383244449476
$ws.Range("M19").Value = 1.047285449029797
$ws.Range("N19").Value = 1.010933693099947
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.013881142900298
$ws.Range("D20").Value = 1.020100588176942
$ws.Range("E20").Value = 1.040665206725609
$ws.Range("F20").Value = 1.04355556857754
$ws.Range("I20").Value = 1.02605289719558
$ws.Range("J20").Value = 1.020229981284428
$ws.Range("K20").Value = 1.023530628463077
$ws.Range("L20").Value = 1.044022529789989
$ws.Range("M20").Value = 1.046902924058897
$ws.Range("N20").Value = 1.010879092302951
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.01293287837518
$ws.Range("D21").Value = 1.019428293046303
$ws.Range("E21").Value = 1.039282058502368
$ws.Range("F21").Value = 1.042101865608308
$ws.Range("I21").Value = 1.025920319867411
$ws.Range("J21").Value = 1.019689217707903
$ws.Range("K21").Value = 1.023071216279645
$ws.Range("L21").Value = 1.042850024779356
$ws.Range("M21").Value = 1.045659442052339
$ws.Range("N21").Value = 1.010701212138899
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.012336497674111
$ws.Range("D22").Value = 1.019004911488335
$ws.Range("E22").Value = 1.038413164830066
$ws.Range("F22").Value = 1.041188542548704
$ws.Range("I22").Value = 1.025834724410998
$ws.Range("J22").Value = 1.019348381968753
$ws.Range("K22").Value = 1.022780931059625
$ws.Range("L22").Value = 1.042112981393368
$ws.Range("M22").Value = 1.044877719690388
$ws.Range("N22").Value = 1.010589091351088
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.01265266610604
$ws.Range("D23").Value = 1.019229418303322
$ws.Range("E23").Value = 1.038873711532943
$ws.Range("F23").Value = 1.04167264863534
$ws.Range("I23").Value = 1.025880312146603
$ws.Range("J23").Value = 1.019529144453989
$ws.Range("K23").Value = 1.022934952377992
$ws.Range("L23").Value = 1.042503687564238
$ws.Range("M23").Value = 1.045292116022376
$ws.Range("N23").Value = 1.010648555198931
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.013897037748154
$ws.Range("D24").Value = 1.020111847607519
$ws.Range("E24").Value = 1.040688408153761
$ws.Range("F24").Value = 1.043579951699181
$ws.Range("I24").Value = 1.026055081863566
$ws.Range("J24").Value = 1.020239033026059
$ws.Range("K24").Value = 1.023538306173571
$ws.Range("L24").Value = 1.044042189747358
$ws.Range("M24").Value = 1.046923773052193
$ws.Range("N24").Value = 1.010882069713413
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.015340295170845
$ws.Range("D25").Value = 1.021132833956144
$ws.Range("E25").Value = 1.042797562077951
$ws.Range("F25").Value = 1.045796270730393
$ws.Range("I25").Value = 1.026248118307597
$ws.Range("J25").Value = 1.021059152354103
$ws.Range("K25").Value = 1.024232182832198
$ws.Range("L25").Value = 1.045828261189635
$ws.Range("M25").Value = 1.04881771244392
$ws.Range("N25").Value = 1.011151820071311
